$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, B (count), C (image), D (word), E (category)
$data = @(
    @(2, 115, "house/house000.jpg", "binden", "house"),
    @(3, 55, "dog/dog011.jpg", "piepen", "dog"),
    @(4, 109, "dog/dog016.jpg", "meinen", "dog"),
    @(5, 49, "dog/dog030.jpg", "kranken", "dog"),
    @(6, 10, "dog/dog017.jpg", "mögen", "dog"),
    @(7, 23, "dog/dog008.jpg", "stoßen", "dog"),
    @(8, 120, "dog/dog021.jpg", "ärgern", "dog"),
    @(9, 16, "house/house008.jpg", "spüren", "house"),
    @(10, 30, "house/house013.jpg", "küssen", "house"),
    @(11, 105, "dog/dog029.jpg", "betteln", "dog"),
    @(12, 0, "house/house014.jpg", "heben", "house"),
    @(13, 113, "house/house007.jpg", "dienen", "house"),
    @(14, 56, "house/house011.jpg", "süßen", "house"),
    @(15, 121, "dog/dog028.jpg", "parken", "dog"),
    @(16, 111, "dog/dog004.jpg", "lernen", "dog"),
    @(17, 119, "house/house031.jpg", "zögern", "house"),
    @(18, 70, "house/house027.jpg", "hacken", "house"),
    @(19, 11, "dog/dog024.jpg", "hassen", "dog"),
    @(20, 74, "dog/dog019.jpg", "hören", "dog"),
    @(21, 34, "dog/dog000.jpg", "passen", "dog"),
    @(22, 75, "house/house009.jpg", "narren", "house"),
    @(23, 6, "house/house010.jpg", "streifen", "house"),
    @(24, 69, "house/house012.jpg", "münzen", "house"),
    @(25, 43, "dog/dog015.jpg", "legen", "dog"),
    @(26, 58, "house/house021.jpg", "duschen", "house"),
    @(27, 50, "dog/dog006.jpg", "grenzen", "dog"),
    @(28, 77, "house/house025.jpg", "spenden", "house"),
    @(29, 100, "dog/dog031.jpg", "nullen", "dog"),
    @(30, 9, "house/house001.jpg", "lügen", "house"),
    @(31, 110, "house/house028.jpg", "wählen", "house"),
    @(32, 32, "house/house020.jpg", "heißen", "house"),
    @(33, 87, "dog/dog009.jpg", "rufen", "dog")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
